# Trade #8 closed at 2026-02-16 22:57:10 - base_strategy DOWN +0.000%
# Append a new trade row (row 9) to both the "All Trades" sheet and the
# "base_strategy" sheet (they mirror each other's trade log).

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(9, 1).Value = 8
    # Leading apostrophe forces text storage so "2026-02-16" isn't
    # reinterpreted as a date serial number (matches column B's existing
    # cells, which are all stored as literal text).
    $ws.Cells.Item(9, 2).Value = "'2026-02-16"
    $ws.Cells.Item(9, 3).Value = "22:57:10"
    $ws.Cells.Item(9, 4).Value = "base_strategy"
    $ws.Cells.Item(9, 5).Value = "DOWN"
    $ws.Cells.Item(9, 6).Value = 0.5
    # Exit Price / Exit Reason are blank text cells for an OPEN trade.
    $ws.Cells.Item(9, 7).Value = ""
    $ws.Cells.Item(9, 8).Value = "OPEN"
    $ws.Cells.Item(9, 9).Value = 0
    $ws.Cells.Item(9, 10).Value = 0
    $ws.Cells.Item(9, 11).Value = 100
    $ws.Cells.Item(9, 12).Value = 0
    $ws.Cells.Item(9, 13).Value = 0
    $ws.Cells.Item(9, 14).Value = 0.6
    $ws.Cells.Item(9, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(9, 16).Value = ""
    $ws.Cells.Item(9, 17).Value = 0
}
